$d = $word.ActiveDocument

$d.Content.Find.Execute("29×49=", $true, $false, $false, $false, $false, $true, 1, $false, "41×91=", 2) | Out-Null
$d.Content.Find.Execute("13×78=", $true, $false, $false, $false, $false, $true, 1, $false, "58×50=", 2) | Out-Null
$d.Content.Find.Execute("95×37=", $true, $false, $false, $false, $false, $true, 1, $false, "28×95=", 2) | Out-Null
$d.Content.Find.Execute("17×37=", $true, $false, $false, $false, $false, $true, 1, $false, "12×12=", 2) | Out-Null
$d.Content.Find.Execute("24×81=", $true, $false, $false, $false, $false, $true, 1, $false, "33×21=", 2) | Out-Null
$d.Content.Find.Execute("80×75=", $true, $false, $false, $false, $false, $true, 1, $false, "31×14=", 2) | Out-Null
$d.Content.Find.Execute("85×21=", $true, $false, $false, $false, $false, $true, 1, $false, "59×69=", 2) | Out-Null
$d.Content.Find.Execute("81×37=", $true, $false, $false, $false, $false, $true, 1, $false, "41×38=", 2) | Out-Null
$d.Content.Find.Execute("75×98=", $true, $false, $false, $false, $false, $true, 1, $false, "74×53=", 2) | Out-Null
$d.Content.Find.Execute("55×19=", $true, $false, $false, $false, $false, $true, 1, $false, "33×59=", 2) | Out-Null
$d.Content.Find.Execute("15×82=", $true, $false, $false, $false, $false, $true, 1, $false, "28×75=", 2) | Out-Null
$d.Content.Find.Execute("21×28=", $true, $false, $false, $false, $false, $true, 1, $false, "43×71=", 2) | Out-Null
$d.Content.Find.Execute("65×24=", $true, $false, $false, $false, $false, $true, 1, $false, "90×30=", 2) | Out-Null
$d.Content.Find.Execute("43×34=", $true, $false, $false, $false, $false, $true, 1, $false, "53×98=", 2) | Out-Null
$d.Content.Find.Execute("39×29=", $true, $false, $false, $false, $false, $true, 1, $false, "11×35=", 2) | Out-Null
$d.Content.Find.Execute("95×29=", $true, $false, $false, $false, $false, $true, 1, $false, "23×22=", 2) | Out-Null
$d.Content.Find.Execute("84×69=", $true, $false, $false, $false, $false, $true, 1, $false, "20×37=", 2) | Out-Null
$d.Content.Find.Execute("17×58=", $true, $false, $false, $false, $false, $true, 1, $false, "82×40=", 2) | Out-Null
$d.Content.Find.Execute("92×20=", $true, $false, $false, $false, $false, $true, 1, $false, "85×84=", 2) | Out-Null
$d.Content.Find.Execute("40×29=", $true, $false, $false, $false, $false, $true, 1, $false, "25×74=", 2) | Out-Null
$d.Content.Find.Execute("25×71=", $true, $false, $false, $false, $false, $true, 1, $false, "63×27=", 2) | Out-Null
$d.Content.Find.Execute("77×83=", $true, $false, $false, $false, $false, $true, 1, $false, "57×87=", 2) | Out-Null
$d.Content.Find.Execute("63×81=", $true, $false, $false, $false, $false, $true, 1, $false, "87×20=", 2) | Out-Null
$d.Content.Find.Execute("87×22=", $true, $false, $false, $false, $false, $true, 1, $false, "72×27=", 2) | Out-Null
$d.Content.Find.Execute("50×16=", $true, $false, $false, $false, $false, $true, 1, $false, "18×56=", 2) | Out-Null
